# Fruta / hortaliza, semanal
# A new weekly price observation is inserted at row 305 (Vega Monumental
# Concepción - Piña, Caramelo, Segunda), pushing all subsequent rows
# (old 305:327) down by one (new 306:328).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 305, shifting existing rows 305:327 -> 306:328.
$ws.Rows(305).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A305").Value = 11
$ws.Range("B305").Value = "Vega Monumental Concepción"
$ws.Range("C305").Value = "Bíobío"
$ws.Range("D305").Value = 45265
$ws.Range("E305").Value = 8
$ws.Range("F305").Value = "Fruta"
$ws.Range("G305").Value = 100108
$ws.Range("H305").Value = "Tropicales y subtropicales"
$ws.Range("I305").Value = 100108005
$ws.Range("J305").Value = "Piña"
$ws.Range("K305").Value = "Caramelo"
$ws.Range("L305").Value = "Segunda"
$ws.Range("M305").Value = 100
$ws.Range("N305").Value = 22000
$ws.Range("O305").Value = 22000
$ws.Range("P305").Value = 22000
$ws.Range("Q305").Value = "$/caja 14 unidades"
$ws.Range("R305").Value = "Ecuador"
$ws.Range("S305").Value = 1571
$ws.Range("T305").Value = 14
